# Generate Report for Handoff
# Update the "Latest Handoff Datetime" (column D) for the 2792364d... file
# row (row 5) on both the zh-cn and de-de status sheets, reflecting a new
# handoff report generation.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D5").Value = "2016-03-07 04:18:27"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D5").Value = "2016-03-07 04:18:37"
